# expenses.xlsx import fixes:
#  - rename the sheet/tab from the old car name to the new one
#  - clear the free-text "Description" entries in column H (rows 2-11)
#    that referenced the old car's specific repair notes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (was "Audi A4 B6 1.8T")
$ws.Name = "Suzuki SX4"

# Clear out the per-row Description text in column H for rows 2-11;
# the sheet keeps the column but the detailed notes are removed/blanked.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 8).Value = " "
}
